{"js": "// Update the dummy test scripts' reported numbers: each of the three\n// per-file detail tables goes from \"1\" test/passed to \"2\" test/passed,\n// and the Summary table's total NumTests/Passed goes from \"3\" to \"6\"\n// (matching the new combined total across the detail tables).\nconst tables = context.document.body.tables;\ntables.load(\"items\");\nawait context.sync();\n\n// Table 0 is the \"Summary\" table (single data row: NumTests, Passed,\n// Failed, Warnings, Errors, Skipped). Bump NumTests and Passed 3 -> 6.\nconst summary = tables.items[0];\nsummary.getCell(1, 0).value = \"6\"; // NumTests\nsummary.getCell(1, 1).value = \"6\"; // Passed\n\n// Tables 1..3 are the per-test-file detail tables (Test, NumTests,\n// Passed, Failed, Warnings, Errors, Skipped). Each has one data row;\n// bump NumTests and Passed 1 -> 2.\nfor (let i = 1; i < tables.items.length; i++) {\n  const detail = tables.items[i];\n  detail.getCell(1, 1).value = \"2\"; // NumTests\n  detail.getCell(1, 2).value = \"2\"; // Passed\n}\n\nawait context.sync();\n", "ps1": "# Update the dummy test scripts' reported numbers: each of the three\n# per-file detail tables goes from \"1\" test/passed to \"2\" test/passed,\n# and the Summary table's total NumTests/Passed goes from \"3\" to \"6\"\n# (matching the new combined total across the detail tables).\n$d = $word.ActiveDocument\n\n# Table 1 is the \"Summary\" table: columns are NumTests, Passed, Failed,\n# Warnings, Errors, Skipped; row 2 is the single data row. Bump\n# NumTests and Passed 3 -> 6.\n$summary = $d.Tables(1)\n$summary.Cell(2, 1).Range.Text = \"6\"  # NumTests\n$summary.Cell(2, 2).Range.Text = \"6\"  # Passed\n\n# Tables 2..4 are the per-test-file detail tables: columns are Test,\n# NumTests, Passed, Failed, Warnings, Errors, Skipped; row 2 is the\n# single data row. Bump NumTests and Passed 1 -> 2.\nfor ($i = 2; $i -le $d.Tables.Count; $i++) {\n    $detail = $d.Tables($i)\n    $detail.Cell(2, 2).Range.Text = \"2\"  # NumTests\n    $detail.Cell(2, 3).Range.Text = \"2\"  # Passed\n}\n"}
